$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Update cell values (text / numbers) for the refreshed
#    correlation-count analysis table (rows 1-12, columns A-D).
# -----------------------------------------------------------------

# Header row
$ws.Cells.Item(1,1).Value = "Factors"
$ws.Cells.Item(1,2).Value = " Correlation Count"
$ws.Cells.Item(1,3).Value = "High Correlation"
$ws.Cells.Item(1,4).Value = "Observation"

# Row 2 - GDP per capita
$ws.Cells.Item(2,1).Value = "GDP per capita"
$ws.Cells.Item(2,2).Value = 11
$ws.Cells.Item(2,3).Value = "Yes"
$ws.Cells.Item(2,4).Value = "Rejecting GDP per capita as a factor of Agricultural production per capita as agricultural production is a subset of GDP."

# Row 3 - Avg temperature
$ws.Cells.Item(3,1).Value = "Avg temperature"
$ws.Cells.Item(3,2).Value = 7
$ws.Cells.Item(3,3).Value = "Yes"
$ws.Cells.Item(3,4).Value = "Avg temperature is a factor beyond a nation's control. "

# Row 4 - Gross enrolment ratio
$ws.Cells.Item(4,1).Value = "Gross enrolment ratio, primary to tertiary, both sexes (%)"
$ws.Cells.Item(4,2).Value = 9
$ws.Cells.Item(4,3).Value = "Yes"
$ws.Cells.Item(4,4).Value = "High-producing nations have high education enrolment ratio. "

# Row 5 - Area
$ws.Cells.Item(5,1).Value = "Area"
$ws.Cells.Item(5,2).Value = 9
$ws.Cells.Item(5,3).Value = "Yes"
$ws.Cells.Item(5,4).Value = "Area available for agriculture is a factor beyond a nation's control. But we can see from plot, even with lower amount of area assigned for agriculture, they have higher agri production per capita."

# Row 6 - Fertilizer Use Per Capita
$ws.Cells.Item(6,1).Value = "Fertilizer Use Per Capita"
$ws.Cells.Item(6,2).Value = 8
$ws.Cells.Item(6,3).Value = "Yes"
$ws.Cells.Item(6,4).Value = "High-producing nations use high amounts of fertilizer per capita. Hence, the governments of low-producing nations can find ways of utilizing more  fertilizer in their agricultural land."

# Row 7 - Water Use Efficiency
$ws.Cells.Item(7,1).Value = "Water Use Efficiency"
$ws.Cells.Item(7,2).Value = 8
$ws.Cells.Item(7,3).Value = "Yes"
$ws.Cells.Item(7,4).Value = "High and low producing nations have same water use efficiency. Need to further analyze which crops are the most water-efficient."

# Row 8 - Credit to Agriculture (new row)
$ws.Cells.Item(8,1).Value = "Credit to Agriculture"
$ws.Cells.Item(8,2).Value = 3
$ws.Cells.Item(8,3).Value = "No"
$ws.Cells.Item(8,4).Value = "Credit to Agriculture may not be a dominant factor for well-producing developed countries but it needs to be further explored as Total agricultural production is  highly correlated with 'Credit to Agriculture' for the entire dataset."

# Row 9 - Agriculture share of Government Expenditure
$ws.Cells.Item(9,1).Value = "Agriculture share of Government Expenditure"
$ws.Cells.Item(9,2).Value = 7
$ws.Cells.Item(9,3).Value = "Yes"
$ws.Cells.Item(9,4).Value = "Even with lower agriculture share of government expenditure, high-producing nations have more production per capita. This does not imply government should spend less on their agriculture sector."

# Row 10 - Gini coefficient
$ws.Cells.Item(10,1).Value = "Gini coefficient"
$ws.Cells.Item(10,2).Value = 6
$ws.Cells.Item(10,3).Value = "Yes"
$ws.Cells.Item(10,4).Value = "Need to be explored further as there is high correlation."

# Row 11 - Population
$ws.Cells.Item(11,1).Value = "Population"
$ws.Cells.Item(11,2).Value = 8
$ws.Cells.Item(11,3).Value = "Yes"
$ws.Cells.Item(11,4).Value = "Rejecting population as a factor as it cannot be easily tunable."

# Row 12 - FDI inflows to Agriculture (new row)
$ws.Cells.Item(12,1).Value = "FDI inflows to Agriculture"
$ws.Cells.Item(12,2).Value = 5
$ws.Cells.Item(12,3).Value = "Yes"
$ws.Cells.Item(12,4).Value = "Need to be explored further as there is high correlation."

# -----------------------------------------------------------------
# 2) Re-apply cell formatting (fill / bold / number-format / wrap)
#    so every cell keeps using the correct existing style.
#    Source cells below already carry the desired style in the
#    original sheet, so copy their formats onto the cells whose
#    style needs to change.
# -----------------------------------------------------------------

$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# 3) Row heights to match the wrapped text of the new content.
# -----------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 29
$ws.Rows.Item(4).RowHeight = 29
$ws.Rows.Item(5).RowHeight = 43.5
$ws.Rows.Item(6).RowHeight = 43.5
$ws.Rows.Item(7).RowHeight = 29
$ws.Rows.Item(8).RowHeight = 58
$ws.Rows.Item(9).RowHeight = 43.5

# -----------------------------------------------------------------
# 4) Widen column D to fit the longer observation text.
# -----------------------------------------------------------------

$ws.Columns.Item(4).ColumnWidth = 83

# -----------------------------------------------------------------
# 5) Restore the last active-cell selection recorded in the file.
# -----------------------------------------------------------------

$null = $ws.Range("D19").Select()
